$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 150
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").Value = -376

$ws.Range("H43").Value = 2658
$ws.Range("I43").Value = 1816.3334
$ws.Range("J43").Value = 3499.6667
$ws.Range("K43").Value = 1816.3334
$ws.Range("L43").Value = 3499.6667
$ws.Range("M43").Value = -1747.3334

$ws.Range("H45").Value = 3109.25
$ws.Range("I45").Value = 518
$ws.Range("J45").Value = 3973
$ws.Range("K45").Value = 1554
$ws.Range("L45").Value = 11919
$ws.Range("M45").Value = -1362
$ws.Range("N45").Value = -12303

$ws.Range("H86").Value = 51770.2
$ws.Range("I86").Value = 57300.223
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 57300.223
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -56177.223
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 51770.2
$ws.Range("I89").Value = 57300.223
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 286501.115
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -280885.115
$ws.Range("N89").Value = -21232

$ws.Range("H113").Value = 2442.0715
$ws.Range("I113").Value = 1711.125
$ws.Range("J113").Value = 3416.6667
$ws.Range("K113").Value = 1711.125
$ws.Range("L113").Value = 3416.6667
$ws.Range("M113").Value = 1542.875

$ws.Range("H116").Value = 26988642
$ws.Range("I116").Value = 33337856
$ws.Range("J116").Value = 4476.5
$ws.Range("K116").Value = 33337856
$ws.Range("L116").Value = 4476.5
$ws.Range("M116").Value = -33334414
$ws.Range("N116").Value = -11360.5

$ws.Range("H125").Value = 6418.4443
$ws.Range("I125").Value = 943.75
$ws.Range("J125").Value = 10798.2
$ws.Range("K125").Value = 8493.75
$ws.Range("L125").Value = 97183.8
$ws.Range("M125").Value = -6033.75

$ws.Range("H137").Value = 24527240
$ws.Range("I137").Value = 1253062.4
$ws.Range("J137").Value = 55559476
$ws.Range("K137").Value = 3759187.2
$ws.Range("L137").Value = 166678428
$ws.Range("M137").Value = -3756637.2
$ws.Range("N137").Value = -166683528


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3503.75
$ws.Range("I61").Value = 1397
$ws.Range("J61").Value = 6759.636
$ws.Range("K61").Value = 1397
$ws.Range("L61").Value = 6759.636
$ws.Range("M61").Value = -1185
$ws.Range("N61").Value = -7183.636

$ws.Range("H74").Value = 19232044
$ws.Range("I74").Value = 20834614
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 20834614
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -20833740
$ws.Range("N74").Value = -2948

$ws.Range("H77").Value = 19232044
$ws.Range("I77").Value = 20834614
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 104173070
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -104168702
$ws.Range("N77").Value = -14736

$ws.Range("H88").Value = 3559.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3559.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3559.8
$ws.Range("M88").ClearContents() | Out-Null
$ws.Range("N88").Value = -4371.8

$ws.Range("H91").Value = 3559.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3559.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3559.8
$ws.Range("M91").ClearContents() | Out-Null
$ws.Range("N91").Value = -6367.8

$ws.Range("H136").Value = 3503.75
$ws.Range("I136").Value = 1397
$ws.Range("J136").Value = 6759.636
$ws.Range("K136").Value = 4191
$ws.Range("L136").Value = 20278.908
$ws.Range("M136").Value = -1641
$ws.Range("N136").Value = -25378.908


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 674.75
$ws.Range("I64").Value = 650
$ws.Range("J64").Value = 683
$ws.Range("K64").Value = 650
$ws.Range("L64").Value = 683
$ws.Range("M64").Value = -425
$ws.Range("N64").Value = -1133

$ws.Range("H67").Value = 674.75
$ws.Range("I67").Value = 650
$ws.Range("J67").Value = 683
$ws.Range("K67").Value = 650
$ws.Range("L67").Value = 683
$ws.Range("M67").Value = 130
$ws.Range("N67").Value = -2243

$ws.Range("H107").Value = 2157.077
$ws.Range("I107").Value = 1929.6666
$ws.Range("J107").Value = 2467.182
$ws.Range("K107").Value = 1929.6666
$ws.Range("L107").Value = 2467.182
$ws.Range("M107").Value = -9.666600000000017


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 86712.8
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 86712.8
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 86712.8
$ws.Range("N9").Value = -87048.8

$ws.Range("H22").Value = 1476.5834
$ws.Range("I22").Value = 1215
$ws.Range("J22").Value = 1999.75
$ws.Range("K22").Value = 1215
$ws.Range("L22").Value = 1999.75
$ws.Range("M22").Value = -865
$ws.Range("N22").Value = -2699.75

$ws.Range("H31").Value = 7175
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 7175
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 7175
$ws.Range("M31").ClearContents() | Out-Null

$ws.Range("H34").Value = 7175
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7175
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7175
$ws.Range("M34").ClearContents() | Out-Null
$ws.Range("N34").Value = -7579

$ws.Range("H58").Value = 1880.7273
$ws.Range("I58").Value = 1828.9
$ws.Range("J58").Value = 2399
$ws.Range("K58").Value = 1828.9
$ws.Range("L58").Value = 2399
$ws.Range("M58").Value = -1625.9

$ws.Range("H105").Value = 2467.8333
$ws.Range("I105").Value = 1935.6666
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1935.6666
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -188.6666

$ws.Range("H136").Value = 1880.7273
$ws.Range("I136").Value = 1828.9
$ws.Range("J136").Value = 2399
$ws.Range("K136").Value = 5486.700000000001
$ws.Range("L136").Value = 7197
$ws.Range("M136").Value = -2936.700000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 143.64285
$ws.Range("I12").Value = 126
$ws.Range("J12").Value = 148.45454
$ws.Range("K12").Value = 378
$ws.Range("L12").Value = 445.36362
$ws.Range("M12").Value = -205
$ws.Range("N12").Value = -791.3636200000001

$ws.Range("H64").Value = 9995
$ws.Range("I64").Value = 6486.25
$ws.Range("J64").Value = 13503.75
$ws.Range("K64").Value = 19458.75
$ws.Range("L64").Value = 40511.25
$ws.Range("M64").Value = -19188.75
$ws.Range("N64").Value = -41051.25

$ws.Range("H67").Value = 9995
$ws.Range("I67").Value = 6486.25
$ws.Range("J67").Value = 13503.75
$ws.Range("K67").Value = 19458.75
$ws.Range("L67").Value = 40511.25
$ws.Range("M67").Value = -18522.75
$ws.Range("N67").Value = -42383.25

$ws.Range("H80").Value = 5014.5713
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5014.5713
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15043.7139
$ws.Range("N80").Value = -16915.7139

$ws.Range("H83").Value = 5014.5713
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5014.5713
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45131.14169999999
$ws.Range("N83").Value = -54491.14169999999

$ws.Range("H107").Value = 1030.875
$ws.Range("I107").Value = 769.8889
$ws.Range("J107").Value = 1187.4667
$ws.Range("K107").Value = 2309.6667
$ws.Range("L107").Value = 3562.4001
$ws.Range("M107").Value = -389.6667000000002

$ws.Range("H112").Value = 3681.8572
$ws.Range("I112").Value = 247.66667
$ws.Range("J112").Value = 6257.5
$ws.Range("K112").Value = 743.00001
$ws.Range("L112").Value = 18772.5
$ws.Range("M112").Value = 364.99999
$ws.Range("N112").Value = -20988.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents() | Out-Null

$ws.Range("H80").Value = 6057.5625
$ws.Range("I80").Value = 2368.4583
$ws.Range("J80").Value = 17124.875
$ws.Range("K80").Value = 2368.4583
$ws.Range("L80").Value = 17124.875
$ws.Range("M80").Value = -1370.4583
$ws.Range("N80").Value = -19120.875

$ws.Range("H83").Value = 6057.5625
$ws.Range("I83").Value = 2368.4583
$ws.Range("J83").Value = 17124.875
$ws.Range("K83").Value = 11842.2915
$ws.Range("L83").Value = 85624.375
$ws.Range("M83").Value = -6850.291499999999
$ws.Range("N83").Value = -95608.375

$ws.Range("H93").Value = 44999.25
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 44999.25
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 44999.25
$ws.Range("N93").Value = -48743.25

$ws.Range("H133").Value = 54772.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 54772.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 54772.332
$ws.Range("N133").Value = -64892.332


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4059.4338
$ws.Range("I7").Value = 3550.8596
$ws.Range("J7").Value = 5174.385
$ws.Range("K7").Value = 3550.8596
$ws.Range("L7").Value = 5174.385
$ws.Range("M7").Value = -3438.8596

$ws.Range("H16").Value = 2709.7058
$ws.Range("I16").Value = 1658.9231
$ws.Range("J16").Value = 6124.75
$ws.Range("K16").Value = 1658.9231
$ws.Range("L16").Value = 6124.75
$ws.Range("M16").Value = -1488.9231

$ws.Range("H61").Value = 4672.154
$ws.Range("I61").Value = 4740.364
$ws.Range("J61").Value = 4297
$ws.Range("K61").Value = 4740.364
$ws.Range("L61").Value = 4297
$ws.Range("M61").Value = -4538.364

$ws.Range("H113").Value = 4672.154
$ws.Range("I113").Value = 4740.364
$ws.Range("J113").Value = 4297
$ws.Range("K113").Value = 4740.364
$ws.Range("L113").Value = 4297
$ws.Range("M113").Value = -2570.364

$ws.Range("H126").Value = 4059.4338
$ws.Range("I126").Value = 3550.8596
$ws.Range("J126").Value = 5174.385
$ws.Range("K126").Value = 10652.5788
$ws.Range("L126").Value = 15523.155
$ws.Range("M126").Value = -8182.578799999999

$ws.Range("H132").Value = 4289.844
$ws.Range("I132").Value = 3680.8809
$ws.Range("J132").Value = 5020.6
$ws.Range("K132").Value = 11042.6427
$ws.Range("L132").Value = 15061.8
$ws.Range("M132").Value = -8512.6427
$ws.Range("N132").Value = -20121.8

$ws.Range("H136").Value = 3881.6924
$ws.Range("I136").Value = 2794.0852
$ws.Range("J136").Value = 6721.5557
$ws.Range("K136").Value = 8382.2556
$ws.Range("L136").Value = 20164.6671
$ws.Range("M136").Value = -5832.2556
$ws.Range("N136").Value = -25264.6671


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3183.65
$ws.Range("I122").Value = 2212.182
$ws.Range("J122").Value = 4371
$ws.Range("K122").Value = 6636.545999999999
$ws.Range("L122").Value = 13113
$ws.Range("M122").Value = -4186.545999999999

$ws.Range("H132").Value = 1614.1464
$ws.Range("I132").Value = 790.5238000000001
$ws.Range("J132").Value = 1897.6885
$ws.Range("K132").Value = 2371.5714
$ws.Range("L132").Value = 5693.0655
$ws.Range("M132").Value = 158.4285999999997

$ws.Range("H136").Value = 3000.682
$ws.Range("I136").Value = 2106.158
$ws.Range("J136").Value = 8666
$ws.Range("K136").Value = 6318.474
$ws.Range("L136").Value = 25998
$ws.Range("M136").Value = -3768.474

